# Renamed Thermdat to Nasa. Changed excel import function. Bug fixes.
#
# Content change in this workbook: the "thermo_model" description (E2) and
# the actual thermo model values used by the two species rows (E3, E4) are
# updated to reflect the renamed/refactored model identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the comment/description of the thermo_model column.
$ws.Range("E2").Value = "Type of thermodynamic model. Supported models include IdealGas and Harmonic"

# Rows 3 & 4: update the actual thermo_model value used for each species
# from the old fully-qualified class path to the new short model name.
$ws.Range("E3").Value = "IdealGas"
$ws.Range("E4").Value = "IdealGas"
